$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: replace with the new listing (生成AI活用 / 日本人モデル画像生成歓迎) ---
$ws.Range("A2").Value = "2025-11-02 06:24:50"
$ws.Range("B2").Value = "【業務自動化×補助金対応】生成AI活用/日本人モデル画像生成歓迎"
$ws.Range("D2").Value = "3,000,000 円 ~ 5,000,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5405834"
$ws.Range("G2").Value = 395
$ws.Range("H2").Value = "🔥AI,Ai ◆自動化"

# --- Row 3: only the capture timestamp changes ---
$ws.Range("A3").Value = "2025-11-02 06:24:50"

# --- Row 4: replace with what used to be row 5's listing (ECフォース EFO) ---
$ws.Range("A4").Value = "2025-11-02 06:24:50"
$ws.Range("B4").Value = "【急募】ECフォース EFO利用 META広告計測設定の経験者募集"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5425263"
$ws.Range("G4").Value = 13

# --- Drop the old rows 5 and 6 (電子工作 listing and the duplicate EC row) ---
$ws.Rows("5:6").Delete()

# --- Column H grows from 12 to 13 characters wide ---
$ws.Columns.Item(8).ColumnWidth = 12.17

# --- Rebuild the hyperlinks collection so it only references F2:F4 with the
#     correct, current targets (stale links to the deleted rows are dropped) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5405834")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5425363")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5425263")
